# "error handling added to savasdialog. Started with export settings to excel workbook"
#
# Functional change captured by this script: add a new "Settings" worksheet
# (placed right after Sheet1) that will hold export/settings key-value pairs,
# starting with a "Use external source" flag (True/False dropdown).

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# --- add the new Settings sheet, right after Sheet1 ---------------------
$settings = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$settings.Name = "Settings"

# Use Excel's standard (US Letter) default page margins for the brand new
# sheet, same as a freshly inserted worksheet would get.
$settings.PageSetup.LeftMargin = 0.75 * 72
$settings.PageSetup.RightMargin = 0.75 * 72
$settings.PageSetup.TopMargin = 1 * 72
$settings.PageSetup.BottomMargin = 1 * 72
$settings.PageSetup.HeaderMargin = 0.5 * 72
$settings.PageSetup.FooterMargin = 0.5 * 72

# --- headers --------------------------------------------------------------
$settings.Range("A1").Value = "Name"
$settings.Range("B1").Value = "Value"

# --- first setting row: Use external source = True ------------------------
$settings.Range("A2").Value = "Use external source"

# Enter the value through a formula + paste-values round trip so Excel
# stores it as plain text ("True") rather than auto-coercing it to a
# native boolean.
$settings.Range("B2").Formula = "=""True"""
$settings.Range("B2").Copy()
$settings.Range("B2").PasteSpecial(-4163)

# Restrict the value to a True/False dropdown list.
$settings.Range("B2").Validation.Add(3, 1, 1, """True,False""")

# --- selection / activation so Settings becomes the active tab ------------
[void]$settings.Range("A5").Select()
[void]$settings.Activate()
